$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 65; this shifts the existing rows 65..115 down to 66..116
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with the new weekly entry
$ws.Cells.Item(65, 1).Value = 11
$ws.Cells.Item(65, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(65, 3).Value = "Bíobío"
$ws.Cells.Item(65, 4).Value = 45167
$ws.Cells.Item(65, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(65, 5).Value = 8
$ws.Cells.Item(65, 6).Value = 100112013
$ws.Cells.Item(65, 7).Value = "Alcachofa"
$ws.Cells.Item(65, 8).Value = "Argentina(o)"
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 120
$ws.Cells.Item(65, 11).Value = 10000
$ws.Cells.Item(65, 12).Value = 10000
$ws.Cells.Item(65, 13).Value = 10000
$ws.Cells.Item(65, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(65, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(65, 16).Value = 200
$ws.Cells.Item(65, 17).Value = 50
$ws.Cells.Item(65, 18).Value = "Hortaliza"
